$wb = $excel.ActiveWorkbook

# Reposition the workbook window (best-effort; some hosts don't persist this).
$win = $wb.Windows.Item(1)
$win.Left = 28680
$win.Top = -120

# Rename sheets: drop the "Soafia_" prefix from each sheet's name.
$wb.Worksheets.Item(1).Name = "GrainBulkheads"
$wb.Worksheets.Item(2).Name = "GrainBulkheadsPlace"
$wb.Worksheets.Item(3).Name = "CargoCompartments"
$wb.Worksheets.Item(4).Name = "CargoCompartmentsParts"
$wb.Worksheets.Item(5).Name = "сontainer_deck_plan"

# Move the active tab / selection from the 5th sheet to the 1st sheet,
# and update the selected cell on the first sheet from C25 to C29.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Select()
$ws1.Range("C29").Select()
